# Updated symbol list with GitHub Actions
# Applies the price/volume/coin-swap updates described by the diff.
#
# The Price column (D) stores numeric-looking values as TEXT in the
# original workbook (t="inlineStr"). A plain $cell.Value = "242.88"
# assignment lets Excel auto-detect the number and store it as a real
# numeric cell, which would change the cell type. To keep these cells
# as text (matching the source data), we enter them with a leading
# apostrophe (Excel's "force text" prefix) and then reset NumberFormat
# back to General so no stray text-format style lingers on the cell.

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.NumberFormat = "General"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Simple price column (D) updates
Set-TextValue $ws "D2"  "242.88"
Set-TextValue $ws "D3"  "23.19"
Set-TextValue $ws "D4"  "5.752"
Set-TextValue $ws "D5"  "0.05801"
Set-TextValue $ws "D6"  "3.419"
Set-TextValue $ws "D7"  "6.472"
Set-TextValue $ws "D9"  "0.8015"
Set-TextValue $ws "D10" "0.1453"
Set-TextValue $ws "D11" "0.07651"
Set-TextValue $ws "D12" "0.03249"
Set-TextValue $ws "D14" "0.09234"
Set-TextValue $ws "D15" "0.001666"
Set-TextValue $ws "D17" "0.04757"

# Row 18: price + volume(1h) text update
Set-TextValue $ws "D18" "0.0005998"
$ws.Range("E18").Value = "17OneONEWorstin24h"

Set-TextValue $ws "D19" "0.006210"
Set-TextValue $ws "D20" "0.005384"
Set-TextValue $ws "D21" "0.001065"
Set-TextValue $ws "D22" "0.0001501"
Set-TextValue $ws "D23" "3.692"
Set-TextValue $ws "D25" "0.3320"
Set-TextValue $ws "D26" "0.1242"
Set-TextValue $ws "D27" "0.0006737"
Set-TextValue $ws "D40" "0.04292"
Set-TextValue $ws "D41" "0.007077"

# Rows 42/43: coins BKEXToken and CEJI swap places (with refreshed prices)
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D42" "0.003603"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D43" "0.1053"
$ws.Range("E43").Value = "42BKEXTokenBKK"

Set-TextValue $ws "D44" "0.009719"

# Row 45: remove trailing "Worstin24h" marker from volume text
$ws.Range("E45").Value = "44ACDXExchangeACXT"

Set-TextValue $ws "D46" "0.00005628"
Set-TextValue $ws "D48" "0.7864"
Set-TextValue $ws "D49" "0.09915"
